$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107, pushing existing rows 107-167 down to 108-168.
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new data point.
$ws.Cells.Item(107, 1).Value = 9
$ws.Cells.Item(107, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(107, 3).Value = "Metropolitana"
$ws.Cells.Item(107, 4).Value = 44488
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 300000001
$ws.Cells.Item(107, 7).Value = "Rabanito"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 7900
$ws.Cells.Item(107, 11).Value = 3000
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = 3494
$ws.Cells.Item(107, 14).Value = '$/cien unidades (volumen en unidades)'
$ws.Cells.Item(107, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(107, 16).Value = 35
$ws.Cells.Item(107, 17).Value = 100
$ws.Cells.Item(107, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D (style index 2, used by D2:D168).
$ws.Cells.Item(107, 4).NumberFormat = $ws.Cells.Item(108, 4).NumberFormat
